$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Insert the new rows -----------------------------------------------
# Final layout (rows 10..19):
#   10 GaúchoDaFronteira F.C   (NEW)
#   11 Esquadrão Gazembrino    (existing, shifted)
#   12  NHU PORÃ SAF.          (existing, shifted)
#   13 SC 100 Sono             (NEW)
#   14 GrioTeam                (NEW)
#   15 GE Bebum                (existing, shifted)
#   16 bugredasmissões         (existing, shifted)
#   17 Pontaç0 F.C.            (NEW)
#   18 lsauer fc               (existing, shifted)
#   19 Grêmio_Campeão_LA_27    (existing, shifted)
#
# Insert from the bottom up (in original row numbering) so each insertion
# point is still valid when we reach it.
$ws.Rows.Item(14).Insert()
$ws.Rows.Item(12).Insert()
$ws.Rows.Item(12).Insert()
$ws.Rows.Item(10).Insert()

# --- 2. Give the new column-A cells the same style as their neighbours ----
# (bold / bordered / centred "s=1" style used throughout column A)
$ws.Range("A9").Copy()
$ws.Range("A10").PasteSpecial(-4122)
$ws.Range("A13").PasteSpecial(-4122)
$ws.Range("A14").PasteSpecial(-4122)
$ws.Range("A17").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- 3. Write the final values for every data row (2..19) -----------------
$ws.Range("A2").Value = 32966
$ws.Range("B2").Value = "La Primeira Patada Es Nuestra"
$ws.Range("C2").Value = 0

$ws.Range("A3").Value = 184499
$ws.Range("B3").Value = "SC ÉoINTER! "
$ws.Range("C3").Value = 0

$ws.Range("A4").Value = 186283
$ws.Range("B4").Value = "FBC Colorado"
$ws.Range("C4").Value = 0

$ws.Range("A5").Value = 287965
$ws.Range("B5").Value = "Doug Leal F.C"
$ws.Range("C5").Value = 0

$ws.Range("A6").Value = 1273719
$ws.Range("B6").Value = "Texas Club 2026"
$ws.Range("C6").Value = 0

$ws.Range("A7").Value = 1488983
$ws.Range("B7").Value = "C R Juvenal"
$ws.Range("C7").Value = 0

$ws.Range("A8").Value = 1747619
$ws.Range("B8").Value = "JV5 Tricolor Gaúcho"
$ws.Range("C8").Value = 0

$ws.Range("A9").Value = 1867254
$ws.Range("B9").Value = "Medonho´s F.C. "
$ws.Range("C9").Value = 0

$ws.Range("A10").Value = 2371918
$ws.Range("B10").Value = "GaúchoDaFronteira F.C"
$ws.Range("C10").Value = 0

$ws.Range("A11").Value = 2916559
$ws.Range("B11").Value = "Esquadrão Gazembrino"
$ws.Range("C11").Value = 0

$ws.Range("A12").Value = 4088673
$ws.Range("B12").Value = " NHU PORÃ SAF."
$ws.Range("C12").Value = 0

$ws.Range("A13").Value = 14709358
$ws.Range("B13").Value = "SC 100 Sono"
$ws.Range("C13").Value = 0

$ws.Range("A14").Value = 14933455
$ws.Range("B14").Value = "GrioTeam"
$ws.Range("C14").Value = 0

$ws.Range("A15").Value = 16411206
$ws.Range("B15").Value = "GE Bebum"
$ws.Range("C15").Value = 0

$ws.Range("A16").Value = 19209079
$ws.Range("B16").Value = "bugredasmissões "
$ws.Range("C16").Value = 0

$ws.Range("A17").Value = 20651178
$ws.Range("B17").Value = "Pontaç0 F.C."
$ws.Range("C17").Value = 0

$ws.Range("A18").Value = 44810918
$ws.Range("B18").Value = "lsauer fc"
$ws.Range("C18").Value = 0

$ws.Range("A19").Value = 47775950
$ws.Range("B19").Value = "Grêmio_Campeão_LA_27"
$ws.Range("C19").Value = 0

Write-Host "Workbook updated: rows 2-19 set, dimension now A1:U19"
